# Auto-generated edit script applying the cryptos.xlsx GitHub Actions update
# (Wed Sep  4 05:30:33 UTC 2024) - updates Price/Volume(1h) columns and
# corrects the Aave/Bittensor row ordering (rows 45 & 47).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.381.48"
$ws.Range("E2").Value = "'  -4.80%  "
$ws.Range("D3").Value = "'2.361.29"
$ws.Range("E3").Value = "'  -6.42%  "
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'512.68"
$ws.Range("E5").Value = "'  -4.61%  "
$ws.Range("D6").Value = "'127.44"
$ws.Range("E6").Value = "'  -5.74%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "'  -0.28%  "
$ws.Range("E8").Value = "'  -2.60%  "
$ws.Range("D9").Value = "'2.378.15"
$ws.Range("E9").Value = "'  -5.72%  "
$ws.Range("D10").Value = "'0.0959"
$ws.Range("E10").Value = "'  -3.78%  "
$ws.Range("E11").Value = "'  -1.25%  "
$ws.Range("E12").Value = "'  -8.17%  "
$ws.Range("D13").Value = "'0.315"
$ws.Range("E13").Value = "'  -5.56%  "
$ws.Range("D14").Value = "'2.781.47"
$ws.Range("E14").Value = "'  -6.23%  "
$ws.Range("D15").Value = "'56.363.92"
$ws.Range("E15").Value = "'  -4.76%  "
$ws.Range("D16").Value = "'21.40"
$ws.Range("E16").Value = "'  -4.77%  "
$ws.Range("E17").Value = "'  -4.81%  "
$ws.Range("D18").Value = "'2.409.33"
$ws.Range("E18").Value = "'  -4.33%  "
$ws.Range("E19").Value = "'  -3.97%  "
$ws.Range("D20").Value = "'4.05"
$ws.Range("E20").Value = "'  -4.79%  "
$ws.Range("D21").Value = "'310.03"
$ws.Range("E21").Value = "'  -3.72%  "
$ws.Range("D22").Value = "'6.17"
$ws.Range("E22").Value = "'  -0.30%  "
$ws.Range("D23").Value = "'0.996"
$ws.Range("E23").Value = "'  -0.45%  "
$ws.Range("D24").Value = "'65.50"
$ws.Range("E24").Value = "'  -0.54%  "
$ws.Range("E25").Value = "'  +0.55%  "
$ws.Range("E26").Value = "'  -5.51%  "
$ws.Range("D27").Value = "'2.468.38"
$ws.Range("E27").Value = "'  -6.04%  "
$ws.Range("E28").Value = "'  -4.86%  "
$ws.Range("D29").Value = "'7.21"
$ws.Range("E29").Value = "'  -4.30%  "
$ws.Range("D30").Value = "'174.83"
$ws.Range("E30").Value = "'  +1.93%  "
$ws.Range("E31").Value = "'  -4.31%  "
$ws.Range("D32").Value = "'0.0₃0714"
$ws.Range("E32").Value = "'  -7.03%  "
$ws.Range("D33").Value = "'6.13"
$ws.Range("E33").Value = "'  -3.12%  "
$ws.Range("D34").Value = "'1.11"
$ws.Range("E34").Value = "'  -7.95%  "
$ws.Range("E35").Value = "'  -0.11%  "
$ws.Range("D36").Value = "'0.996"
$ws.Range("E36").Value = "'  -0.37%  "
$ws.Range("D37").Value = "'17.65"
$ws.Range("E37").Value = "'  -3.08%  "
$ws.Range("E38").Value = "'  -5.82%  "
$ws.Range("E39").Value = "'  -7.32%  "
$ws.Range("D40").Value = "'0.803"
$ws.Range("E40").Value = "'  +0.88%  "
$ws.Range("D41").Value = "'35.72"
$ws.Range("E41").Value = "'  -2.65%  "
$ws.Range("E42").Value = "'  -6.97%  "
$ws.Range("E43").Value = "'  -3.91%  "
$ws.Range("D44").Value = "'4.86"
$ws.Range("E44").Value = "'  -4.77%  "
$ws.Range("B45").Value = "'Aave"
$ws.Range("C45").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'123.14"
$ws.Range("E45").Value = "'  -6.97%  "
$ws.Range("D46").Value = "'0.568"
$ws.Range("E46").Value = "'  -4.65%  "
$ws.Range("B47").Value = "'Bittensor"
$ws.Range("C47").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "'252.27"
$ws.Range("E47").Value = "'  -9.77%  "
$ws.Range("E48").Value = "'  -3.30%  "
$ws.Range("D49").Value = "'0.0488"
$ws.Range("E49").Value = "'  -4.77%  "
$ws.Range("D50").Value = "'0.0207"
$ws.Range("E50").Value = "'  -5.82%  "
$ws.Range("D51").Value = "'16.64"
$ws.Range("E51").Value = "'  -6.70%  "
